# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.429.76'
$ws.Range('E2').Value = '  +1.60%  '
$ws.Range('D3').Value = '2.280.47'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.97'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.91'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +6.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.530'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.97'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +11.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0801'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.71'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = '2.633.51'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.47'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').Value = '2.286.10'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.799'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.36%  '
$ws.Range('D18').Value = '42.288.24'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.59'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '0.0₃0913'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.72'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.94'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.60'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.89'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.93'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +9.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.56'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '160.59'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.25'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  +3.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0746'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.18'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.19%  '
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('E39').Value = '  +3.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.15'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.39'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +14.32%  '
$ws.Range('D43').Value = '2.002.33'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0287'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.06'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.06'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.40'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.53'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.45'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '92.77'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.77%  '
